$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H10").Value = 3756.6
$ws1.Range("I10").Value = 1260
$ws1.Range("H12").Value = 1128.6
$ws1.Range("Q18").Value = 2215.7
$ws1.Range("H22").Value = "2 de 20"
$ws1.Range("I22").Value = "1 de 20"
$ws1.Range("Q22").Value = "2 de 20"

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = 4687.09
$ws2.Range("F12").Value = 4357.27
$ws2.Range("F18").Value = 2215.7
$ws2.Range("F22").Value = 12201.87

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D7").Value = 4885.2
$ws3.Range("E7").Value = -3085.2
$ws3.Range("F7").Value = 2.714

$ws3.Range("D8").Value = 1260
$ws3.Range("E8").Value = -635
$ws3.Range("F8").Value = 2.016

$ws3.Range("D14").Value = 2647.34
$ws3.Range("E14").Value = -2164.34
$ws3.Range("F14").Value = 5.481035196687371

$ws3.Range("D19").Value = 12201.87
$ws3.Range("E19").Value = 38185.32762291768
$ws3.Range("F19").Value = 0.2421621081472927
